$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells that would otherwise be parsed as numbers
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '98.021.11'
$ws.Range("E2").Value = '  -0.46%  '

$ws.Range("D3").Value = '3.400.89'
$ws.Range("E3").Value = '  +1.00%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").Value = '253.90'
$ws.Range("E5").Value = '  -1.64%  '

$ws.Range("D6").Value = '678.05'
$ws.Range("E6").Value = '  +1.76%  '

$ws.Range("D7").Value = '1.44'
$ws.Range("E7").Value = '  -6.60%  '

$ws.Range("D8").Value = '0.430'
$ws.Range("E8").Value = '  -7.55%  '

$ws.Range("E9").Value = '  -3.06%  '

$ws.Range("E10").Value = '  +0.01%  '

$ws.Range("D11").Value = '3.397.45'
$ws.Range("E11").Value = '  +0.98%  '

$ws.Range("E12").Value = '  +0.89%  '

$ws.Range("D13").Value = '41.56'
$ws.Range("E13").Value = '  -1.51%  '

$ws.Range("D14").Value = '6.26'
$ws.Range("E14").Value = '  +10.74%  '

$ws.Range("D15").Value = '97.784.57'
$ws.Range("E15").Value = '  +0.84%  '

$ws.Range("D16").Value = '0.0000264'
$ws.Range("E16").Value = '  -2.40%  '

$ws.Range("D17").Value = '4.041.36'
$ws.Range("E17").Value = '  +1.22%  '

$ws.Range("D18").Value = '8.91'
$ws.Range("E18").Value = '  +16.48%  '

$ws.Range("D19").Value = '3.388.34'
$ws.Range("E19").Value = '  +0.62%  '

$ws.Range("D20").Value = '0.567'
$ws.Range("E20").Value = '  +29.60%  '

$ws.Range("D21").Value = '17.34'
$ws.Range("E21").Value = '  +2.64%  '

$ws.Range("D22").Value = '11.01'
$ws.Range("E22").Value = '  +4.43%  '

$ws.Range("D24").Value = '506.89'
$ws.Range("E24").Value = '  -4.64%  '

$ws.Range("E25").Value = '  -7.34%  '

$ws.Range("D26").Value = '6.53'
$ws.Range("E26").Value = '  +4.42%  '

$ws.Range("D27").Value = '99.21'
$ws.Range("E27").Value = '  -3.49%  '

$ws.Range("D28").Value = '12.56'
$ws.Range("E28").Value = '  -0.60%  '

$ws.Range("D29").Value = '3.579.61'
$ws.Range("E29").Value = '  +0.92%  '

$ws.Range("E30").Value = '  -0.16%  '

$ws.Range("D31").Value = '11.49'
$ws.Range("E31").Value = '  +3.89%  '

$ws.Range("E32").Value = '  +0.24%  '

$ws.Range("E33").Value = '  +1.69%  '

$ws.Range("D34").Value = '2.59'
$ws.Range("E34").Value = '  +21.82%  '

$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.05%  '

$ws.Range("D36").Value = '0.567'
$ws.Range("E36").Value = '  +2.76%  '

$ws.Range("D37").Value = '29.38'
$ws.Range("E37").Value = '  -0.62%  '

$ws.Range("D38").Value = '1.51'
$ws.Range("E38").Value = '  +11.40%  '

$ws.Range("D39").Value = '7.90'
$ws.Range("E39").Value = '  +0.11%  '

$ws.Range("D40").Value = '532.03'
$ws.Range("E40").Value = '  +0.48%  '

$ws.Range("E41").Value = '  -3.62%  '

$ws.Range("D43").Value = '24.72'

$ws.Range("D44").Value = '0.867'
$ws.Range("E44").Value = '  +2.98%  '

$ws.Range("B45").Value = 'Cosmos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D45").Value = '8.95'
$ws.Range("E45").Value = '  +11.99%  '

$ws.Range("D46").Value = '0.0431'
$ws.Range("E46").Value = '  -2.11%  '

$ws.Range("B47").Value = 'MantraDAO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D47").Value = '3.75'
$ws.Range("E47").Value = '  -0.26%  '

$ws.Range("D48").Value = '1.73'
$ws.Range("E48").Value = '  +12.95%  '

$ws.Range("D49").Value = '5.73'
$ws.Range("E49").Value = '  +11.55%  '

$ws.Range("D50").Value = '55.58'
$ws.Range("E50").Value = '  +10.65%  '

$ws.Range("D51").Value = '3.20'
$ws.Range("E51").Value = '  -6.91%  '
